# #327 Ajout des profils d'acces a58d18c1e8091c98efec92c8c093b361a253eee5
#
# 1. Bump the "Date" metadata value.
# 2. Swap the two mapping columns (AK / AL) on the "Elements" sheet:
#    "Mapping: RIM Mapping" and "Mapping: Spécification métier vers
#    l'extension ROR CommuneCog" trade places (header text, the three
#    populated data cells, and the column widths all move together).

$wb = $excel.ActiveWorkbook

# --- 1. Metadata!B8 : Date -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(8, 2).Value2 = "2024-03-19T13:17:15+00:00"

# --- 2. Elements sheet: swap columns AK (37) and AL (38) -------------------
$ws = $wb.Worksheets.Item("Elements")

$colLeft = 37   # AK
$colRight = 38  # AL

# Only touch the rows whose AK/AL pair actually differs (row 2 and row 4 are
# blank on both sides, so leave those untouched rather than rewriting an
# empty-shared-string cell into a plain blank cell).
$rowsToSwap = @(1, 3, 5, 6)
foreach ($r in $rowsToSwap) {
    $leftCell = $ws.Cells.Item($r, $colLeft)
    $rightCell = $ws.Cells.Item($r, $colRight)
    $leftValue = $leftCell.Value2
    $rightValue = $rightCell.Value2
    $leftCell.Value2 = $rightValue
    $rightCell.Value2 = $leftValue
}

# Swap the column widths too, so the wider "Spécification métier" column
# keeps its width after moving to AK and vice versa.
$ws.Columns.Item($colLeft).ColumnWidth = 68.4
$ws.Columns.Item($colRight).ColumnWidth = 24.17
